# Append the next stanza of "Jabberwocky" after "And stood awhile in
# thought." — matching the target diff:
#   - the existing last paragraph gets a second run containing a single
#     trailing space
#   - a blank NoSpacing paragraph follows
#   - four new NoSpacing paragraphs carry the new stanza text (with the
#     same proofErr spell-check wrapping style already used elsewhere in
#     the doc)
#   - the _GoBack bookmark relocates to the very end of the new content,
#     i.e. around the last new paragraph, which is how Word leaves it
#     after the most recent edit lands there

$d = $word.ActiveDocument

$lastPara = $d.Paragraphs.Last
$startPos = $lastPara.Range.Start
$endPos   = $lastPara.Range.End

# Re-seat the paragraph's bounds as a plain Range (not the live
# Paragraph.Range object) so replacing its contents doesn't leave a
# leftover empty paragraph behind for the old paragraph mark.
$rng = $d.Range($startPos, $endPos)

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t xml:space="preserve">  And stood awhile in thought.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr></w:p>
<w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t xml:space="preserve">And, as in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>uffish</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> thought he stood,</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t xml:space="preserve">  The Jabberwock, with eyes of flame,</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t xml:space="preserve">Came whiffling through the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>tulgey</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> wood,</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t xml:space="preserve">  And burbled as it came!</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$rng.InsertXML($xml)
